# Remove unused master slides (slide layouts 17, 18, 19: "Title and
# Content", "Title Slide", "1_Title and Content") that are not referenced
# by any slide in the deck.
$p = $ppt.ActivePresentation
$layouts = $p.SlideMaster.CustomLayouts

# Delete from the end so earlier indices are not shifted while removing.
$layouts.Item(19).Delete()
$layouts.Item(18).Delete()
$layouts.Item(17).Delete()
